# "added another test case"
#
# - Rename the "busniess flow" sheet to "test_busniess_flow".
# - Move the active/selected tab from "test_register" to "test_busniess_flow".
# - Update each sheet's remembered selection (active cell):
#     test_register:       J10 -> F4
#     test_busniess_flow:  F3  -> D13

$wb = $excel.ActiveWorkbook

$wsRegister = $wb.Worksheets.Item("test_register")
$wsFlow     = $wb.Worksheets.Item("busniess flow")

# Rename the third sheet.
$wsFlow.Name = "test_busniess_flow"

# test_register is no longer the selected tab; just move its remembered
# selection to F4.
$wsRegister.Activate()
$wsRegister.Range("F4").Select()

# test_busniess_flow becomes the active tab, with D13 selected.
$wsFlow.Activate()
$wsFlow.Range("D13").Select()
